# Rewinding to the last functional status until weathermap is integrated
#
# Updates the GHI_2024-02-03.xlsx workbook:
#  - sunrise/sunset timestamps (shared across "Daily" and "Hourly" sheets)
#  - longitude ("lon", column B) on every data row of both sheets
#  - Daily-sheet row-2 clear/cloudy sky aggregates
#  - Hourly-sheet per-hour clear/cloudy sky values for hours 7-17 (rows 9-19)

$wb = $excel.ActiveWorkbook

$wsDaily = $wb.Worksheets.Item("Daily")
$wsHourly = $wb.Worksheets.Item("Hourly")

$newSunrise = "2024-02-03T07:41:36"
$newSunset  = "2024-02-03T17:28:35"
$newLon     = 24.724419

# --- sunrise / sunset (column E / F) -------------------------------------
# Touch every cell that references the old shared strings so the engine
# can recycle/replace the shared-string-table entries instead of leaving
# stale, unreferenced duplicates behind.
$wsDaily.Range("E2").Value = $newSunrise
$wsDaily.Range("F2").Value = $newSunset

for ($r = 2; $r -le 25; $r++) {
    $wsHourly.Cells.Item($r, 5).Value = $newSunrise
    $wsHourly.Cells.Item($r, 6).Value = $newSunset
}

# --- longitude (column B) on every data row -------------------------------
$wsDaily.Range("B2").Value = $newLon

for ($r = 2; $r -le 25; $r++) {
    $wsHourly.Cells.Item($r, 2).Value = $newLon
}

# --- Daily sheet row 2: clear/cloudy sky aggregates -----------------------
$wsDaily.Range("G2").Value = 2574.57
$wsDaily.Range("H2").Value = 5760.33
$wsDaily.Range("I2").Value = 673.95
$wsDaily.Range("J2").Value = 650.8
$wsDaily.Range("K2").Value = 0
$wsDaily.Range("L2").Value = 650.8

# --- Hourly sheet: per-hour clear/cloudy sky values -----------------------
# columns: G=hour, H=clear_sky_ghi, I=clear_sky_dni, J=clear_sky_dhi,
#          K=cloudy_sky_ghi, L=cloudy_sky_dni, M=cloudy_sky_dhi

$hourly = @{
    9  = @{ H = 1.31;               I = 12.17;  J = 2.91;               K = 0.45;               M = 0.45 }
    10 = @{ H = 76.04000000000001;  I = 353.18; J = 41.06;              K = 22.36;              M = 22.36 }
    11 = @{ H = 211.65;             I = 596.48; J = 67.45999999999999;  K = 52.91;              M = 52.91 }
    12 = @{ H = 330.42;             I = 706.28; J = 82.08;              K = 82.59999999999999;  M = 82.59999999999999 }
    13 = @{ H = 410.41;             I = 759.49; J = 90.01000000000001;  K = 102.6;  L = 0; M = 102.6 }
    14 = @{ H = 441.23;             I = 777.17; J = 92.79000000000001;  K = 110.31; L = 0; M = 110.31 }
    15 = @{ H = 419.27;             I = 764.77; J = 90.81;              K = 104.9;  L = 0; M = 104.9 }
    16 = @{ H = 347.08;             I = 718.49; J = 83.81;              K = 87.56;  L = 0; M = 87.56 }
    17 = @{ H = 233.81;             I = 621.1;  J = 70.52;              K = 60.17;  L = 0; M = 60.17 }
    18 = @{ H = 98.44;              I = 412.44; J = 46.92;              K = 25.67;  L = 0; M = 25.67 }
    19 = @{ H = 4.91;               I = 38.76;  J = 5.57;               K = 1.27;              M = 1.27 }
}

foreach ($row in $hourly.Keys) {
    $vals = $hourly[$row]
    foreach ($col in $vals.Keys) {
        $wsHourly.Range("$col$row").Value = $vals[$col]
    }
}
